$wb = $excel.ActiveWorkbook

# New date being appended as row 64 on every price sheet.
$newDate = "'2025-05-04"

# Sheet name -> new price value (carried forward from row 63 / 2025-05-03).
$updates = @{
    "N-Dense"                  = "'38"
    "N-Type"                   = "'37.3"
    "N-type Wafer"              = "'1.02"
    "Cell Topcon 183mm"         = "'0.273"
    "Module Topcon 183mm"       = "'0.09"
    "Silver Rear_side"          = "'5,360"
    "Silver Busbar front-side"  = "'8,025"
    "Silver finger front-side"  = "'8,075"
    "USD_CNY"                   = "'7.2927"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A64").Value = $newDate
    $ws.Range("B64").Value = $updates[$sheetName]

    # Keep the cells on the default "Normal" style so no stray number
    # formatting (e.g. quote-prefix / date formatting) gets attached.
    $ws.Range("A64").Style = "Normal"
    $ws.Range("B64").Style = "Normal"
}
